$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1949.5
$ws.Range("I19").Value = 1882.8334
$ws.Range("J19").Value = 2149.5
$ws.Range("K19").Value = 1882.8334
$ws.Range("L19").Value = 2149.5
$ws.Range("M19").Value = -1707.8334
$ws.Range("N19").Value = -2499.5

$ws.Range("H21").Value = 3922
$ws.Range("I21").Value = 3922
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3922
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -3454

$ws.Range("H23").Value = 3922
$ws.Range("I23").Value = 3922
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 3922
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3688

$ws.Range("H40").Value = 5172.6665
$ws.Range("I40").Value = 5164.6665
$ws.Range("J40").Value = 5180.6665
$ws.Range("K40").Value = 5164.6665
$ws.Range("L40").Value = 5180.6665
$ws.Range("M40").Value = -4989.6665
$ws.Range("N40").Value = -5530.6665

$ws.Range("H62").Value = 4249
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4249
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4249
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5497

$ws.Range("H65").Value = 4249
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4249
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 21245
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -27485

$ws.Range("H88").Value = 1369.0476
$ws.Range("I88").Value = 1690.9
$ws.Range("J88").Value = 1076.4546
$ws.Range("K88").Value = 1690.9
$ws.Range("L88").Value = 1076.4546
$ws.Range("M88").Value = -1284.9
$ws.Range("N88").Value = -1888.4546

$ws.Range("H91").Value = 1369.0476
$ws.Range("I91").Value = 1690.9
$ws.Range("J91").Value = 1076.4546
$ws.Range("K91").Value = 1690.9
$ws.Range("L91").Value = 1076.4546
$ws.Range("M91").Value = -286.9000000000001
$ws.Range("N91").Value = -3884.4546

$ws.Range("H106").Value = 5358.7646
$ws.Range("I106").Value = 5073.3335
$ws.Range("J106").Value = 7499.5
$ws.Range("K106").Value = 5073.3335
$ws.Range("L106").Value = 7499.5
$ws.Range("M106").Value = -4442.3335
$ws.Range("N106").Value = -8761.5

$ws.Range("H112").Value = 5250
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5250
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 15750
$ws.Range("N112").Value = -17966

$ws.Range("H113").Value = 6477.5
$ws.Range("I113").Value = 6477.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6477.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3223.5

$ws.Range("H132").Value = 4731.3335
$ws.Range("I132").Value = 5544.5
$ws.Range("J132").Value = 665.5
$ws.Range("K132").Value = 16633.5
$ws.Range("L132").Value = 1996.5
$ws.Range("M132").Value = -14103.5
$ws.Range("N132").Value = -7056.5

$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1950
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 5750
$ws.Range("I138").Value = 5750
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 17250
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -12110

$ws.Range("H141").Value = 1669.04
$ws.Range("I141").Value = 1718.4584
$ws.Range("J141").Value = 483
$ws.Range("K141").Value = 5155.3752
$ws.Range("L141").Value = 1449
$ws.Range("M141").Value = 24.6247999999996
$ws.Range("N141").Value = -11809

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1998
$ws.Range("I45").Value = 1550.3334
$ws.Range("J45").Value = 2893.3333
$ws.Range("K45").Value = 1550.3334
$ws.Range("L45").Value = 2893.3333
$ws.Range("M45").Value = -1173.3334
$ws.Range("N45").Value = -3647.3333

$ws.Range("H61").Value = 5459.4
$ws.Range("I61").Value = 920.4
$ws.Range("J61").Value = 9998.4
$ws.Range("K61").Value = 920.4
$ws.Range("L61").Value = 9998.4
$ws.Range("M61").Value = -708.4
$ws.Range("N61").Value = -10422.4

$ws.Range("H122").Value = 2010.3334
$ws.Range("I122").Value = 2010.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6031.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3581.0002

$ws.Range("H136").Value = 5459.4
$ws.Range("I136").Value = 920.4
$ws.Range("J136").Value = 9998.4
$ws.Range("K136").Value = 2761.2
$ws.Range("L136").Value = 29995.2
$ws.Range("M136").Value = -211.1999999999998
$ws.Range("N136").Value = -35095.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7958.391
$ws.Range("I86").Value = 5577.7856
$ws.Range("J86").Value = 11661.556
$ws.Range("K86").Value = 5577.7856
$ws.Range("L86").Value = 11661.556
$ws.Range("M86").Value = -4454.7856
$ws.Range("N86").Value = -13907.556

$ws.Range("H89").Value = 7958.391
$ws.Range("I89").Value = 5577.7856
$ws.Range("J89").Value = 11661.556
$ws.Range("K89").Value = 27888.928
$ws.Range("L89").Value = 58307.78
$ws.Range("M89").Value = -22272.928
$ws.Range("N89").Value = -69539.78

$ws.Range("H96").Value = 20738.572
$ws.Range("I96").Value = 20738.572
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 20738.572
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -17992.572

$ws.Range("H134").Value = 5729.6333
$ws.Range("I134").Value = 5739.963
$ws.Range("J134").Value = 5636.6665
$ws.Range("K134").Value = 17219.889
$ws.Range("L134").Value = 16909.9995
$ws.Range("M134").Value = -14684.889
$ws.Range("N134").Value = -21979.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7122.125
$ws.Range("I31").Value = 4396
$ws.Range("J31").Value = 11665.667
$ws.Range("K31").Value = 4396
$ws.Range("L31").Value = 11665.667
$ws.Range("M31").Value = -4101
$ws.Range("N31").Value = -12255.667

$ws.Range("H34").Value = 7122.125
$ws.Range("I34").Value = 4396
$ws.Range("J34").Value = 11665.667
$ws.Range("K34").Value = 4396
$ws.Range("L34").Value = 11665.667
$ws.Range("M34").Value = -4194
$ws.Range("N34").Value = -12069.667

$ws.Range("H105").Value = 1723.6666
$ws.Range("I105").Value = 1788.4
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 1788.4
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -41.40000000000009
$ws.Range("N105").Value = -4894

$ws.Range("H134").Value = 2460.6316
$ws.Range("I134").Value = 2298.647
$ws.Range("J134").Value = 3837.5
$ws.Range("K134").Value = 6895.941
$ws.Range("L134").Value = 11512.5
$ws.Range("M134").Value = -4360.941
$ws.Range("N134").Value = -16582.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2907
$ws.Range("I5").Value = 2907
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 8721
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -8609
$ws.Range("N5").ClearContents()

$ws.Range("H135").Value = 2907
$ws.Range("I135").Value = 2907
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 26163
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -23628
$ws.Range("N135").ClearContents()

$ws.Range("H139").Value = 7531
$ws.Range("I139").Value = 7531
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 22593
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -17453

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2967
$ws.Range("I80").Value = 3032.3333
$ws.Range("J80").Value = 2901.6667
$ws.Range("K80").Value = 3032.3333
$ws.Range("L80").Value = 2901.6667
$ws.Range("M80").Value = -2034.3333
$ws.Range("N80").Value = -4897.6667

$ws.Range("H83").Value = 2967
$ws.Range("I83").Value = 3032.3333
$ws.Range("J83").Value = 2901.6667
$ws.Range("K83").Value = 15161.6665
$ws.Range("L83").Value = 14508.3335
$ws.Range("M83").Value = -10169.6665
$ws.Range("N83").Value = -24492.3335

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 342.54544
$ws.Range("I55").Value = 231.28572
$ws.Range("J55").Value = 537.25
$ws.Range("K55").Value = 231.28572
$ws.Range("L55").Value = 537.25
$ws.Range("M55").Value = -58.28572
$ws.Range("N55").Value = -883.25

$ws.Range("H132").Value = 2259
$ws.Range("I132").Value = 2393.9167
$ws.Range("J132").Value = 640
$ws.Range("K132").Value = 7181.750100000001
$ws.Range("L132").Value = 1920
$ws.Range("M132").Value = -4651.750100000001
$ws.Range("N132").Value = -6980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 885.7778
$ws.Range("I107").Value = 946.125
$ws.Range("J107").Value = 403
$ws.Range("K107").Value = 2838.375
$ws.Range("L107").Value = 1209
$ws.Range("M107").Value = -918.375
$ws.Range("N107").Value = -5049

$ws.Range("H132").Value = 5299.6665
$ws.Range("I132").Value = 4671
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 14013
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -11483
$ws.Range("N132").Value = -27560

$ws.Range("H136").Value = 6552.4614
$ws.Range("I136").Value = 4658.6
$ws.Range("J136").Value = 53899
$ws.Range("K136").Value = 13975.8
$ws.Range("L136").Value = 161697
$ws.Range("M136").Value = -11425.8
$ws.Range("N136").Value = -166797

Write-Output "Applied all Phantom_Profits updates"